$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.742.57'
$ws.Range("E2").Value = '  +0.13%  '

$ws.Range("D3").Value = '2.831.20'
$ws.Range("E3").Value = '  +2.53%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = '''352.85'
$ws.Range("E5").Value = '  +5.76%  '

$ws.Range("D6").Value = '''113.15'
$ws.Range("E6").Value = '  -2.71%  '

$ws.Range("E7").Value = '  +3.96%  '

$ws.Range("D8").Value = '''0.999'
$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("D9").Value = '''0.599'
$ws.Range("E9").Value = '  +3.97%  '

$ws.Range("D10").Value = '''41.53'
$ws.Range("E10").Value = '  -0.89%  '

$ws.Range("E11").Value = '  -1.07%  '

$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = '''0.131'
$ws.Range("E12").Value = '  +1.41%  '

$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").Value = '''19.92'
$ws.Range("E13").Value = '  -1.20%  '

$ws.Range("D14").Value = '''7.71'
$ws.Range("E14").Value = '  +0.58%  '

$ws.Range("D15").Value = '3.278.52'
$ws.Range("E15").Value = '  +2.58%  '

$ws.Range("D16").Value = '2.830.82'
$ws.Range("E16").Value = '  +1.69%  '

$ws.Range("D17").Value = '''0.894'
$ws.Range("E17").Value = '  +0.48%  '

$ws.Range("D18").Value = '51.628.28'
$ws.Range("E18").Value = '  -0.06%  '

$ws.Range("E19").Value = '  +7.41%  '

$ws.Range("E20").Value = '  -1.92%  '

$ws.Range("D21").Value = '''13.44'
$ws.Range("E21").Value = '  -0.75%  '

$ws.Range("D22").Value = '0.0₃0992'
$ws.Range("E22").Value = '  +1.92%  '

$ws.Range("D23").Value = '''270.50'
$ws.Range("E23").Value = '  -2.87%  '

$ws.Range("D24").Value = '''69.66'
$ws.Range("E24").Value = '  +0.06%  '

$ws.Range("D25").Value = '''2.76'
$ws.Range("E25").Value = '  +3.02%  '

$ws.Range("D26").Value = '''26.68'
$ws.Range("E26").Value = '  -0.52%  '

$ws.Range("E27").Value = '  +0.07%  '

$ws.Range("D28").Value = '''10.27'
$ws.Range("E28").Value = '  +0.85%  '

$ws.Range("E29").Value = '  +1.25%  '

$ws.Range("E30").Value = '  -1.26%  '

$ws.Range("D31").Value = '''50.66'
$ws.Range("E31").Value = '  +1.25%  '

$ws.Range("D32").Value = '''33.88'
$ws.Range("E32").Value = '  -3.47%  '

$ws.Range("D33").Value = '''0.0446'
$ws.Range("E33").Value = '  +26.70%  '

$ws.Range("D34").Value = '''5.80'
$ws.Range("E34").Value = '  +4.11%  '

$ws.Range("D35").Value = '''0.0825'
$ws.Range("E35").Value = '  -0.08%  '

$ws.Range("E36").Value = '  -0.07%  '

$ws.Range("D38").Value = '''3.21'
$ws.Range("E38").Value = '  -0.75%  '

$ws.Range("D39").Value = '''4.87'
$ws.Range("E39").Value = '  -2.86%  '

$ws.Range("D40").Value = '''18.02'
$ws.Range("E40").Value = '  -5.05%  '

$ws.Range("D41").Value = '''23.50'
$ws.Range("E41").Value = '  +1.06%  '

$ws.Range("E42").Value = '  +4.33%  '

$ws.Range("D43").Value = '''0.116'
$ws.Range("E43").Value = '  +1.05%  '

$ws.Range("D44").Value = '''125.41'

$ws.Range("E45").Value = '  +0.25%  '

$ws.Range("D46").Value = '2.079.07'
$ws.Range("E46").Value = '  -0.57%  '

$ws.Range("D47").Value = '''3.31'
$ws.Range("E47").Value = '  -0.14%  '

$ws.Range("E48").Value = '  +3.65%  '

$ws.Range("E49").Value = '  +3.06%  '

$ws.Range("E50").Value = '  +6.28%  '

$ws.Range("D51").Value = '''60.73'
$ws.Range("E51").Value = '  +1.36%  '

